$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New teams for two additional leagues (Ligue 1 + Belgian Pro League),
# appended below the existing Ligue 2 / Bundesliga data and sorted by long name.
$newTeams = @(
    ,@("Angers", "Angers SCO")
    ,@("AS Monaco", "AS Monaco FC")
    ,@("St Etienne", "AS Saint-Etienne")
    ,@("Dijon FCO", "Dijon FCO")
    ,@("Bordeaux", "FC Girondins Bordeaux")
    ,@("Lorient", "FC Lorient")
    ,@("Metz", "FC Metz")
    ,@("Nantes", "FC Nantes")
    ,@("Lille", "LOSC Lille")
    ,@("Montpellier", "Montpellier HSC")
    ,@("Nice", "OGC Nice")
    ,@("Lyon", "Olympique Lyonnais")
    ,@("Marseille", "Olympique Marseille")
    ,@("Nimes", "Olympique Nimes")
    ,@("Paris Saint-Germain", "Paris Saint-Germain FC")
    ,@("Lens", "RC Lens")
    ,@("Strasbourg", "RC Strasbourg Alsace")
    ,@("Brest", "Stade Brestois 29")
    ,@("Reims", "Stade Reims")
    ,@("Stade Rennes", "Stade Rennais FC")
    ,@("Cercle Brugge", "Cercle Brugge")
    ,@("Club Brugge", "Club Brugge")
    ,@("KFCO Beerschot-Wilrijk", "K Beerschot VA")
    ,@("KAA Gent", "KAA Gent")
    ,@("Eupen", "KAS Eupen")
    ,@("Genk", "KRC Genk")
    ,@("KV Kortrijk", "KV Kortrijk")
    ,@("KV Oostende", "KV Oostende")
    ,@("OH Leuven", "Oud-Heverlee Leuven")
    ,@("Antwerp", "Royal Antwerp FC")
    ,@("Sporting de Charleroi", "Royal Charleroi SC")
    ,@("Mouscron-Peruwelz", "Royal Excel Mouscron")
    ,@("Anderlecht", "RSC Anderlecht")
    ,@("St. Truidense", "St. Truidense VV")
    ,@("Standard Liege", "Standard Liege")
    ,@("SV Zulte Waregem", "SV Zulte Waregem")
    ,@("Waasland-Beveren", "Waasland-Beveren")
    ,@("KV Mechelen", "Yellow-Red KV Mechelen")
)

$startRow = 40
for ($i = 0; $i -lt $newTeams.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value2 = $newTeams[$i][0]
    $ws.Cells.Item($row, 2).Value2 = $newTeams[$i][1]
}

# Re-apply the sort Excel remembers for the table so the sortState
# metadata covers the freshly added range as well.
$ws.Range("B22:B39").Sort($ws.Range("B22"))
$ws.Range("B40:B59").Sort($ws.Range("B40"))
$ws.Range("B60:B77").Sort($ws.Range("B60"))

# Match the scroll position / selection left behind in the saved file.
$win = $excel.ActiveWindow
$win.ScrollRow = 54
$win.ScrollColumn = 1
$ws.Range("W66").Select()

